# Database.xlsx - "Added passwords for every account"
#
# UserList sheet (sheet1): every account row gets a password hash in column H
# and a flag (1/0) in column G. The previously-empty row 5 is removed so the
# trailing test accounts (e, f, g, banana) shift up one row and the old last
# row (9) collapses into row 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The blank row 5 disappears; rows 6-9 shift up to become rows 5-8.
$ws.Rows.Item(5).Delete()

# Column G: password-assigned flag.
$ws.Range("G1").Value = 1
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 0

# Column H: password hashes for each account.
$ws.Range("H1").Value = "5e884898da28047151d0e56f8dc6292773603d0d6aabbdd62a11ef721d1542d8"
$ws.Range("H2").Value = "b14d501a594442a01c6859541bcb3e8164d183d32937b851835442f69d5c94e"
$ws.Range("H3").Value = "6cf615d5bcaac778352a8f1f3360d23f02f34ec182e259897fd6ce485d7870d4"
$ws.Range("H4").Value = "5906ac361a137e2d286465cd6588ebb5ac3f5ae955001100bc41577c3d751764"
$ws.Range("H5").Value = "b97873a40f73abedd8d685a7cd5e5f85e4a9cfb83eac26886640a0813850122b"
$ws.Range("H6").Value = "8b2c86ea9cf2ea4eb517fd1e06b74f399e7fec0fef92e3b482a6cf2e2b092023"
$ws.Range("H7").Value = "598a1a400c1dfdf36974e69d7e1bc98593f2e15015eed8e9b7e47a83b31693d5"
$ws.Range("H8").Value = "5860836e8f13fc9837539a597d4086bfc0299e54ad92148d54538b5c3feefb7c"

# Reflect the last-selected cell when the workbook was saved.
$ws.Range("G4").Select()
